# Insert a new data row above current row 29 (shifts existing rows 29-75 down to 30-76)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new record's data
$ws.Range("A29").Value = 7
$ws.Range("B29").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C29").Value = "Ñuble"
$ws.Range("D29").Value = 45259
$ws.Range("E29").Value = 16
$ws.Range("F29").Value = 300000000
$ws.Range("G29").Value = "Espárragos"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 200
$ws.Range("K29").Value = 1500
$ws.Range("L29").Value = 1500
$ws.Range("M29").Value = 1500
$ws.Range("N29").Value = '$/kilo'
$ws.Range("O29").Value = "Provincia de Diguillín"
$ws.Range("P29").Value = 1500
$ws.Range("Q29").Value = 1
$ws.Range("R29").Value = "Hortaliza"
